$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 393; this shifts the existing rows 393-428
# down to 394-429 and keeps their data intact.
$ws.Rows.Item(393).Insert()

# Populate the newly inserted row 393 with the new weekly price record.
$ws.Range("A393").Value2 = 4
$ws.Range("B393").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C393").Value2 = "Los Lagos"
$ws.Range("D393").Value2 = 45106
$ws.Range("D393").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E393").Value2 = 10
$ws.Range("F393").Value2 = 100112021
$ws.Range("G393").Value2 = "Ají"
$ws.Range("H393").Value2 = "Inferno"
$ws.Range("I393").Value2 = "Primera"
$ws.Range("J393").Value2 = 70
$ws.Range("K393").Value2 = 22000
$ws.Range("L393").Value2 = 22000
$ws.Range("M393").Value2 = 22000
$ws.Range("N393").Value2 = "`$/caja 10 kilos"
$ws.Range("O393").Value2 = "Región de Arica y Parinacota"
$ws.Range("P393").Value2 = 2200
$ws.Range("Q393").Value2 = 10
$ws.Range("R393").Value2 = "Hortaliza"
